$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.615.23"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.631.21"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.83"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.492"
$ws.Range("E7").Value = "  +1.17%  "
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.94"
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("E11").Value = "  +3.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.859.09"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.632.98"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.622.15"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.96"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.13"
$ws.Range("E19").Value = "  +3.90%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.27"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.17"
$ws.Range("E23").Value = "  +2.32%  "
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.33"
$ws.Range("E25").Value = "  +2.14%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.86"
$ws.Range("E28").Value = "  +4.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.36"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0518"
$ws.Range("E30").Value = "  +3.55%  "
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.164.17"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.808"
$ws.Range("E38").Value = "  +2.41%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.32"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.790"
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.771.39"
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.42"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.43"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0513"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.410"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.52"
$ws.Range("E50").Value = "  +4.15%  "
$ws.Range("E51").Value = "  +0.05%  "
